# Trade #22 closed at 2026-02-17 04:08:27 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: Total Trades (B6) and Win Rate % (B9)
#  - Strategy Status sheet: MarketMaking row Trades (D4) and Win Rate % (G4)
#  - All Trades sheet: append new trade row 23
#  - MarketMaking sheet: append new trade row 23 (mirrors All Trades)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 22
$summary.Range("B9").Value = 31.82

# ---------------------------------------------------------------------------
# 2. Strategy Status sheet (MarketMaking row is row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 22
$status.Range("G4").Value = 31.82

# ---------------------------------------------------------------------------
# 3. Helper to write the new trade row (row 23) into a trades-style sheet
# ---------------------------------------------------------------------------
function Add-TradeRow22([object]$ws) {
    $ws.Range("A23").Value = 22

    # Date/time columns are stored as plain text in this workbook, not as
    # Excel date/time values, so force them to remain text.
    $ws.Range("B23").Value = "'2026-02-17"
    $ws.Range("C23").Value = "04:08:21"

    $ws.Range("D23").Value = "MarketMaking"
    $ws.Range("E23").Value = "DOWN"
    $ws.Range("F23").Value = 0.61
    $ws.Range("G23").Value = 0.61
    $ws.Range("H23").Value = "CLOSED"
    $ws.Range("I23").Value = 0
    $ws.Range("J23").Value = 0
    $ws.Range("K23").Value = 100.02
    $ws.Range("L23").Value = 0
    $ws.Range("M23").Value = 0
    $ws.Range("N23").Value = 0.6
    $ws.Range("O23").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P23").Value = "early_exit"
    $ws.Range("Q23").Value = 0.12
}

# ---------------------------------------------------------------------------
# 4. All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow22 $allTrades

# ---------------------------------------------------------------------------
# 5. MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow22 $marketMaking
